$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from the last existing
# header cell (AC1) onto the three new header cells before setting values.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the team record (Wins/Losses/Ties) for every data row (2-55).
$ws.Range("AD2:AD55").Value = 86
$ws.Range("AE2:AE55").Value = 76
$ws.Range("AF2:AF55").Value = 0
